$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("manually_curated_test_set")

# Rows 87, 88, 89 become "blank template" rows: the RANDBETWEEN-driven
# helper columns (L/N) were erroring out (#NUM!) because the paragraph
# had not actually been coded yet, so every helper column D..P is reset
# back to the same "not yet coded" placeholder pattern used by the other
# still-empty rows (e.g. row 86).
$blankRows = 87, 88, 89
foreach ($r in $blankRows) {
    $ws.Range("D$r").Value = $false
    $ws.Range("E$r").Value = "NA"
    $ws.Range("F$r").Value = $false
    $ws.Range("G$r").Value = "NA"
    $ws.Range("H$r").Value = "NA"
    $ws.Range("I$r").Value = $false
    $ws.Range("J$r").Value = "NA"
    $ws.Range("K$r").Value = 0
    $ws.Range("L$r").Value = "NA"
    $ws.Range("M$r").Value = "NA"
    $ws.Range("N$r").Value = "NA"
    $ws.Range("O$r").Value = "NA"
    $ws.Range("P$r").Value = "NA"
}

# Row 90 is now fully curated with real data.
$ws.Range("D90").Value = $false
$ws.Range("E90").Value = "NA"
$ws.Range("F90").Value = $false
$ws.Range("G90").Value = 139
$ws.Range("H90").Value = $false
$ws.Range("I90").Value = $true
$ws.Range("J90").Value = "The Cronbach's alpha for internal consistency reliability was r = 0.76, p < 0.05 suggesting a reasonable reliability."
$ws.Range("K90").Value = 4
$ws.Range("L90").Value = 2
$ws.Range("M90").Value = 3
$ws.Range("N90").Value = 3
$ws.Range("O90").Value = "Df 2,218 F = 10.77, p = .000"
$ws.Range("P90").Value = "Results"

# Leave the scroll position on the newly curated row.
[void]$ws.Range("D90").Select()
